$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C; this shifts the existing FlowRate_GPM column (C) to D.
$ws.Columns("C").Insert()

# New header for inserted column C, and the new stddev column E.
$ws.Range("C1").Value = "TapPressure_H20_stddev"
$ws.Range("E1").Value = "FlowRate_GPM_stddev"

# Column E isn't covered by the C-insert shift, so match the workbook's
# center-aligned cell style used throughout the rest of the sheet.
$ws.Range("E1:E10").HorizontalAlignment = -4108

# TapPressure_H20_stddev values (column C)
$ws.Range("C2").Value = 0.61808575456808545
$ws.Range("C3").Value = 0.95223421488623394
$ws.Range("C4").Value = 0.58579006478430684
$ws.Range("C5").Value = 0.90483700189591587
$ws.Range("C6").Value = 6.8159848884808998
$ws.Range("C7").Value = 1.9994824330311083
$ws.Range("C8").Value = 0.50707987536481469
$ws.Range("C9").Value = 1.0803101406540592
$ws.Range("C10").Value = 0.95276964687169663

# FlowRate_GPM_stddev values (column E)
$ws.Range("E2").Value = 0.095812316536027828
$ws.Range("E3").Value = 0.075033325929214939
$ws.Range("E4").Value = 0.043243496620878424
$ws.Range("E5").Value = 0.052630789467762389
$ws.Range("E6").Value = 0.087349871207690624
$ws.Range("E7").Value = 0.14300349646075278
$ws.Range("E8").Value = 0.044944410108488694
$ws.Range("E9").Value = 0.088204308284797778
$ws.Range("E10").Value = 0.10616025621672183

# Column widths to match target layout (nearest value reachable through
# the ColumnWidth setter's internal rounding).
$ws.Columns("C").ColumnWidth = 21
$ws.Columns("D").ColumnWidth = 12.8333333333
$ws.Columns("E").ColumnWidth = 19.3333333333

# Update selection to mirror the saved worksheet view.
$ws.Range("E15").Select()
